# redmine #9271 - cal sheet corrections for GA05MOAS-GL494
# 1) "Moorings"!J2 (Cruise Number): "AT26-30" -> "AT-26-30", with the
#    inserted dash highlighted in blue, rest of the text left as-is.
# 2) "Asset_Cal_Info"!F4 (CC_angular_resolution calibration value):
#    1.13 -> 1.096, re-keyed in a distinguishing blue Calibri font.

$wb = $excel.ActiveWorkbook

$moorings = $wb.Worksheets.Item("Moorings")
$cruise = $moorings.Range("J2")
$cruise.Value = "AT-26-30"

# Color just the inserted "-" character blue (matches the rich-text run
# split produced when a single character is re-colored in Excel).
$dash = $cruise.Characters(3, 1)
$dash.Font.Color = 16711680

$assetCal = $wb.Worksheets.Item("Asset_Cal_Info")
$coeff = $assetCal.Range("F4")
$coeff.Value = 1.096
$coeff.Font.Name = "Calibri"
$coeff.Font.Size = 10
$coeff.Font.Color = 16711680

# Restore the selections left behind in the saved workbook (Moorings stays
# the active/front tab, as it was originally).
$assetCal.Activate()
$assetCal.Range("F4").Select()

$moorings.Activate()
$moorings.Range("J3").Select()
$excel.ActiveWindow.Zoom = 100
